$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'327.88"
$ws.Range("E2").Formula = "'-0.92%"
$ws.Range("D3").Formula = "'43.89"
$ws.Range("E3").Formula = "'5.47%"
$ws.Range("D4").Formula = "'5.408"
$ws.Range("E4").Formula = "'-5.10%"
$ws.Range("D5").Formula = "'0.08085"
$ws.Range("E5").Formula = "'-3.69%"
$ws.Range("D6").Formula = "'8.702"
$ws.Range("E6").Formula = "'-1.27%"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Formula = "'4.306"
$ws.Range("E7").Formula = "'-3.77%"
$ws.Range("B8").Value = "FTXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D8").Formula = "'1.900"
$ws.Range("E8").Formula = "'-5.49%"
$ws.Range("D9").Formula = "'2.723"
$ws.Range("E9").Formula = "'-6.12%"
$ws.Range("D10").Formula = "'0.9413"
$ws.Range("E10").Formula = "'1.87%"
$ws.Range("E11").Formula = "'-5.09%"
$ws.Range("D12").Formula = "'0.1894"
$ws.Range("E12").Formula = "'-4.45%"
$ws.Range("D13").Formula = "'0.09511"
$ws.Range("E13").Formula = "'0.33%"
$ws.Range("D14").Formula = "'0.04144"
$ws.Range("E14").Formula = "'7.67%"
$ws.Range("D15").Formula = "'0.1071"
$ws.Range("E15").Formula = "'0.91%"
$ws.Range("D16").Formula = "'0.001274"
$ws.Range("E16").Formula = "'-2.24%"
$ws.Range("D17").Formula = "'0.006082"
$ws.Range("E17").Formula = "'-0.41%"
$ws.Range("D18").Formula = "'3.575"
$ws.Range("E18").Formula = "'4.37%"
$ws.Range("D20").Formula = "'8.505"
$ws.Range("E20").Formula = "'-3.07%"
$ws.Range("D21").Formula = "'0.1359"
$ws.Range("E21").Formula = "'-0.26%"
$ws.Range("D22").Formula = "'0.2606"
$ws.Range("E22").Formula = "'3.85%"
$ws.Range("D23").Formula = "'0.04371"
$ws.Range("E23").Formula = "'-1.02%"
$ws.Range("E24").Formula = "'-2.64%"
$ws.Range("D25").Formula = "'0.004291"
$ws.Range("E25").Formula = "'-2.41%"
$ws.Range("D26").Formula = "'0.0001235"
$ws.Range("E26").Formula = "'1.21%"
$ws.Range("D27").Formula = "'0.0004018"
$ws.Range("E27").Formula = "'0.67%"
$ws.Range("D39").Formula = "'0.02665"
$ws.Range("E39").Formula = "'-6.52%"
$ws.Range("D40").Formula = "'0.05417"
$ws.Range("E40").Formula = "'-2.05%"
$ws.Range("D41").Formula = "'0.007691"
$ws.Range("E41").Formula = "'-3.36%"
$ws.Range("D42").Formula = "'0.009769"
$ws.Range("E42").Formula = "'8.49%"
$ws.Range("D43").Formula = "'0.1392"
$ws.Range("E43").Formula = "'-2.85%"
$ws.Range("D44").Formula = "'0.002128"
$ws.Range("E44").Formula = "'2.74%"
$ws.Range("D45").Formula = "'0.009896"
$ws.Range("E45").Formula = "'-15.60%"
$ws.Range("D46").Formula = "'0.00007082"
$ws.Range("E46").Formula = "'1.95%"
$ws.Range("D47").Formula = "'0.00000000753"
$ws.Range("E47").Formula = "'0.38%"
$ws.Range("D48").Formula = "'0.003549"
$ws.Range("E48").Formula = "'2.36%"
$ws.Range("D49").Formula = "'0.002286"
$ws.Range("E49").Formula = "'0.33%"
$ws.Range("D50").Formula = "'0.00002109"
$ws.Range("E50").Formula = "'0.38%"
$ws.Range("D51").Formula = "'0.0002009"
$ws.Range("E51").Formula = "'0.38%"
